$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Formula = "=""67.378.51"""
$ws.Cells.Item(2,4).Copy()
$ws.Cells.Item(2,4).PasteSpecial(-4163)
$ws.Cells.Item(2,5).Formula = "=""  -0.06%  """
$ws.Cells.Item(2,5).Copy()
$ws.Cells.Item(2,5).PasteSpecial(-4163)
$ws.Cells.Item(3,4).Formula = "=""2.620.76"""
$ws.Cells.Item(3,4).Copy()
$ws.Cells.Item(3,4).PasteSpecial(-4163)
$ws.Cells.Item(3,5).Formula = "=""  -1.75%  """
$ws.Cells.Item(3,5).Copy()
$ws.Cells.Item(3,5).PasteSpecial(-4163)
$ws.Cells.Item(4,5).Formula = "=""  -0.03%  """
$ws.Cells.Item(4,5).Copy()
$ws.Cells.Item(4,5).PasteSpecial(-4163)
$ws.Cells.Item(5,4).Formula = "=""593.75"""
$ws.Cells.Item(5,4).Copy()
$ws.Cells.Item(5,4).PasteSpecial(-4163)
$ws.Cells.Item(5,5).Formula = "=""  -0.93%  """
$ws.Cells.Item(5,5).Copy()
$ws.Cells.Item(5,5).PasteSpecial(-4163)
$ws.Cells.Item(6,4).Formula = "=""167.41"""
$ws.Cells.Item(6,4).Copy()
$ws.Cells.Item(6,4).PasteSpecial(-4163)
$ws.Cells.Item(6,5).Formula = "=""  +0.60%  """
$ws.Cells.Item(6,5).Copy()
$ws.Cells.Item(6,5).PasteSpecial(-4163)
$ws.Cells.Item(7,5).Formula = "=""  -0.01%  """
$ws.Cells.Item(7,5).Copy()
$ws.Cells.Item(7,5).PasteSpecial(-4163)
$ws.Cells.Item(8,5).Formula = "=""  -2.04%  """
$ws.Cells.Item(8,5).Copy()
$ws.Cells.Item(8,5).PasteSpecial(-4163)
$ws.Cells.Item(9,4).Formula = "=""2.620.26"""
$ws.Cells.Item(9,4).Copy()
$ws.Cells.Item(9,4).PasteSpecial(-4163)
$ws.Cells.Item(9,5).Formula = "=""  -1.75%  """
$ws.Cells.Item(9,5).Copy()
$ws.Cells.Item(9,5).PasteSpecial(-4163)
$ws.Cells.Item(10,5).Formula = "=""  -2.92%  """
$ws.Cells.Item(10,5).Copy()
$ws.Cells.Item(10,5).PasteSpecial(-4163)
$ws.Cells.Item(11,5).Formula = "=""  +1.17%  """
$ws.Cells.Item(11,5).Copy()
$ws.Cells.Item(11,5).PasteSpecial(-4163)
$ws.Cells.Item(12,5).Formula = "=""  +0.82%  """
$ws.Cells.Item(12,5).Copy()
$ws.Cells.Item(12,5).PasteSpecial(-4163)
$ws.Cells.Item(13,5).Formula = "=""  +0.45%  """
$ws.Cells.Item(13,5).Copy()
$ws.Cells.Item(13,5).PasteSpecial(-4163)
$ws.Cells.Item(14,4).Formula = "=""27.62"""
$ws.Cells.Item(14,4).Copy()
$ws.Cells.Item(14,4).PasteSpecial(-4163)
$ws.Cells.Item(14,5).Formula = "=""  -0.56%  """
$ws.Cells.Item(14,5).Copy()
$ws.Cells.Item(14,5).PasteSpecial(-4163)
$ws.Cells.Item(15,4).Formula = "=""3.117.56"""
$ws.Cells.Item(15,4).Copy()
$ws.Cells.Item(15,4).PasteSpecial(-4163)
$ws.Cells.Item(16,5).Formula = "=""  -1.35%  """
$ws.Cells.Item(16,5).Copy()
$ws.Cells.Item(16,5).PasteSpecial(-4163)
$ws.Cells.Item(17,4).Formula = "=""67.505.06"""
$ws.Cells.Item(17,4).Copy()
$ws.Cells.Item(17,4).PasteSpecial(-4163)
$ws.Cells.Item(17,5).Formula = "=""  +0.10%  """
$ws.Cells.Item(17,5).Copy()
$ws.Cells.Item(17,5).PasteSpecial(-4163)
$ws.Cells.Item(18,4).Formula = "=""2.627.87"""
$ws.Cells.Item(18,4).Copy()
$ws.Cells.Item(18,4).PasteSpecial(-4163)
$ws.Cells.Item(18,5).Formula = "=""  -1.79%  """
$ws.Cells.Item(18,5).Copy()
$ws.Cells.Item(18,5).PasteSpecial(-4163)
$ws.Cells.Item(19,4).Formula = "=""11.97"""
$ws.Cells.Item(19,4).Copy()
$ws.Cells.Item(19,4).PasteSpecial(-4163)
$ws.Cells.Item(19,5).Formula = "=""  +2.19%  """
$ws.Cells.Item(19,5).Copy()
$ws.Cells.Item(19,5).PasteSpecial(-4163)
$ws.Cells.Item(20,4).Formula = "=""8.00"""
$ws.Cells.Item(20,4).Copy()
$ws.Cells.Item(20,4).PasteSpecial(-4163)
$ws.Cells.Item(20,5).Formula = "=""  +2.53%  """
$ws.Cells.Item(20,5).Copy()
$ws.Cells.Item(20,5).PasteSpecial(-4163)
$ws.Cells.Item(21,4).Formula = "=""357.29"""
$ws.Cells.Item(21,4).Copy()
$ws.Cells.Item(21,4).PasteSpecial(-4163)
$ws.Cells.Item(21,5).Formula = "=""  -1.79%  """
$ws.Cells.Item(21,5).Copy()
$ws.Cells.Item(21,5).PasteSpecial(-4163)
$ws.Cells.Item(22,5).Formula = "=""  -1.26%  """
$ws.Cells.Item(22,5).Copy()
$ws.Cells.Item(22,5).PasteSpecial(-4163)
$ws.Cells.Item(23,5).Formula = "=""  -2.19%  """
$ws.Cells.Item(23,5).Copy()
$ws.Cells.Item(23,5).PasteSpecial(-4163)
$ws.Cells.Item(24,2).Formula = "=""Aptos"""
$ws.Cells.Item(24,2).Copy()
$ws.Cells.Item(24,2).PasteSpecial(-4163)
$ws.Cells.Item(24,3).Formula = "=""https://coinranking.com/coin/HGYj5JCv5+aptos-apt"""
$ws.Cells.Item(24,3).Copy()
$ws.Cells.Item(24,3).PasteSpecial(-4163)
$ws.Cells.Item(24,4).Formula = "=""10.37"""
$ws.Cells.Item(24,4).Copy()
$ws.Cells.Item(24,4).PasteSpecial(-4163)
$ws.Cells.Item(24,5).Formula = "=""  +2.53%  """
$ws.Cells.Item(24,5).Copy()
$ws.Cells.Item(24,5).PasteSpecial(-4163)
$ws.Cells.Item(25,5).Formula = "=""  -0.06%  """
$ws.Cells.Item(25,5).Copy()
$ws.Cells.Item(25,5).PasteSpecial(-4163)
$ws.Cells.Item(26,2).Formula = "=""SuiNetwork"""
$ws.Cells.Item(26,2).Copy()
$ws.Cells.Item(26,2).PasteSpecial(-4163)
$ws.Cells.Item(26,3).Formula = "=""https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"""
$ws.Cells.Item(26,3).Copy()
$ws.Cells.Item(26,3).PasteSpecial(-4163)
$ws.Cells.Item(26,4).Formula = "=""1.94"""
$ws.Cells.Item(26,4).Copy()
$ws.Cells.Item(26,4).PasteSpecial(-4163)
$ws.Cells.Item(26,5).Formula = "=""  -3.90%  """
$ws.Cells.Item(26,5).Copy()
$ws.Cells.Item(26,5).PasteSpecial(-4163)
$ws.Cells.Item(27,4).Formula = "=""69.89"""
$ws.Cells.Item(27,4).Copy()
$ws.Cells.Item(27,4).PasteSpecial(-4163)
$ws.Cells.Item(27,5).Formula = "=""  -1.14%  """
$ws.Cells.Item(27,5).Copy()
$ws.Cells.Item(27,5).PasteSpecial(-4163)
$ws.Cells.Item(28,4).Formula = "=""2.762.86"""
$ws.Cells.Item(28,4).Copy()
$ws.Cells.Item(28,4).PasteSpecial(-4163)
$ws.Cells.Item(29,4).Formula = "=""0.999"""
$ws.Cells.Item(29,4).Copy()
$ws.Cells.Item(29,4).PasteSpecial(-4163)
$ws.Cells.Item(29,5).Formula = "=""  -0.08%  """
$ws.Cells.Item(29,5).Copy()
$ws.Cells.Item(29,5).PasteSpecial(-4163)
$ws.Cells.Item(30,5).Formula = "=""  -1.27%  """
$ws.Cells.Item(30,5).Copy()
$ws.Cells.Item(30,5).PasteSpecial(-4163)
$ws.Cells.Item(31,4).Formula = "=""543.73"""
$ws.Cells.Item(31,4).Copy()
$ws.Cells.Item(31,4).PasteSpecial(-4163)
$ws.Cells.Item(31,5).Formula = "=""  -0.82%  """
$ws.Cells.Item(31,5).Copy()
$ws.Cells.Item(31,5).PasteSpecial(-4163)
$ws.Cells.Item(32,5).Formula = "=""  -0.36%  """
$ws.Cells.Item(32,5).Copy()
$ws.Cells.Item(32,5).PasteSpecial(-4163)
$ws.Cells.Item(33,5).Formula = "=""  -2.37%  """
$ws.Cells.Item(33,5).Copy()
$ws.Cells.Item(33,5).PasteSpecial(-4163)
$ws.Cells.Item(34,5).Formula = "=""  -1.32%  """
$ws.Cells.Item(34,5).Copy()
$ws.Cells.Item(34,5).PasteSpecial(-4163)
$ws.Cells.Item(35,5).Formula = "=""  +4.88%  """
$ws.Cells.Item(35,5).Copy()
$ws.Cells.Item(35,5).PasteSpecial(-4163)
$ws.Cells.Item(37,4).Formula = "=""1.50"""
$ws.Cells.Item(37,4).Copy()
$ws.Cells.Item(37,4).PasteSpecial(-4163)
$ws.Cells.Item(37,5).Formula = "=""  -2.60%  """
$ws.Cells.Item(37,5).Copy()
$ws.Cells.Item(37,5).PasteSpecial(-4163)
$ws.Cells.Item(38,4).Formula = "=""157.72"""
$ws.Cells.Item(38,4).Copy()
$ws.Cells.Item(38,4).PasteSpecial(-4163)
$ws.Cells.Item(38,5).Formula = "=""  +1.46%  """
$ws.Cells.Item(38,5).Copy()
$ws.Cells.Item(38,5).PasteSpecial(-4163)
$ws.Cells.Item(39,4).Formula = "=""19.05"""
$ws.Cells.Item(39,4).Copy()
$ws.Cells.Item(39,4).PasteSpecial(-4163)
$ws.Cells.Item(39,5).Formula = "=""  -1.88%  """
$ws.Cells.Item(39,5).Copy()
$ws.Cells.Item(39,5).PasteSpecial(-4163)
$ws.Cells.Item(40,4).Formula = "=""0.367"""
$ws.Cells.Item(40,4).Copy()
$ws.Cells.Item(40,4).PasteSpecial(-4163)
$ws.Cells.Item(40,5).Formula = "=""  -1.19%  """
$ws.Cells.Item(40,5).Copy()
$ws.Cells.Item(40,5).PasteSpecial(-4163)
$ws.Cells.Item(41,5).Formula = "=""  +1.97%  """
$ws.Cells.Item(41,5).Copy()
$ws.Cells.Item(41,5).PasteSpecial(-4163)
$ws.Cells.Item(42,5).Formula = "=""  -0.94%  """
$ws.Cells.Item(42,5).Copy()
$ws.Cells.Item(42,5).PasteSpecial(-4163)
$ws.Cells.Item(43,5).Formula = "=""  -0.71%  """
$ws.Cells.Item(43,5).Copy()
$ws.Cells.Item(43,5).PasteSpecial(-4163)
$ws.Cells.Item(44,5).Formula = "=""  +0.08%  """
$ws.Cells.Item(44,5).Copy()
$ws.Cells.Item(44,5).PasteSpecial(-4163)
$ws.Cells.Item(45,5).Formula = "=""  -3.27%  """
$ws.Cells.Item(45,5).Copy()
$ws.Cells.Item(45,5).PasteSpecial(-4163)
$ws.Cells.Item(46,4).Formula = "=""0.0₆0300"""
$ws.Cells.Item(46,4).Copy()
$ws.Cells.Item(46,4).PasteSpecial(-4163)
$ws.Cells.Item(46,5).Formula = "=""  +0.16%  """
$ws.Cells.Item(46,5).Copy()
$ws.Cells.Item(46,5).PasteSpecial(-4163)
$ws.Cells.Item(47,4).Formula = "=""152.56"""
$ws.Cells.Item(47,4).Copy()
$ws.Cells.Item(47,4).PasteSpecial(-4163)
$ws.Cells.Item(47,5).Formula = "=""  -0.54%  """
$ws.Cells.Item(47,5).Copy()
$ws.Cells.Item(47,5).PasteSpecial(-4163)
$ws.Cells.Item(48,4).Formula = "=""0.579"""
$ws.Cells.Item(48,4).Copy()
$ws.Cells.Item(48,4).PasteSpecial(-4163)
$ws.Cells.Item(48,5).Formula = "=""  -1.56%  """
$ws.Cells.Item(48,5).Copy()
$ws.Cells.Item(48,5).PasteSpecial(-4163)
$ws.Cells.Item(49,5).Formula = "=""  -1.37%  """
$ws.Cells.Item(49,5).Copy()
$ws.Cells.Item(49,5).PasteSpecial(-4163)
$ws.Cells.Item(50,5).Formula = "=""  -1.09%  """
$ws.Cells.Item(50,5).Copy()
$ws.Cells.Item(50,5).PasteSpecial(-4163)
$ws.Cells.Item(51,4).Formula = "=""0.0771"""
$ws.Cells.Item(51,4).Copy()
$ws.Cells.Item(51,4).PasteSpecial(-4163)
$ws.Cells.Item(51,5).Formula = "=""  -0.47%  """
$ws.Cells.Item(51,5).Copy()
$ws.Cells.Item(51,5).PasteSpecial(-4163)
$excel.CutCopyMode = $false
